$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header columns
$ws.Range("A1").Value = 'mx_state'
$ws.Range("B1").Value = 'mx_municipality'
$ws.Range("C1").Value = 'n_matriculas'
$ws.Range("D1").Value = 'pct_matriculas'
$ws.Range("B5").Value = 'Amatenango De La Frontera'
$ws.Range("A25").Value = 'Ciudad De México'
$ws.Range("B29").Value = 'Cuajimalpa De Morelos'
$ws.Range("A42").Value = 'Estado De México'
$ws.Range("B43").Value = 'Ecatepec De Morelos'
$ws.Range("B45").Value = 'San Felipe Del Progreso'
$ws.Range("B49").Value = 'Tlalnepantla De Baz'
$ws.Range("B52").Value = 'Valle De Bravo'
$ws.Range("B54").Value = 'Apaseo El Alto'
$ws.Range("B57").Value = 'Jaral Del Progreso'
$ws.Range("B63").Value = 'Acapulco De Juárez'
$ws.Range("B65").Value = 'Ajuchitlán Del Progreso'
$ws.Range("B67").Value = 'Atenango Del Río'
$ws.Range("B69").Value = 'Ayutla De Los Libres'
$ws.Range("B70").Value = 'Chilapa De Álvarez'
$ws.Range("B71").Value = 'Chilpancingo De Los Bravo'
$ws.Range("B74").Value = 'Coyuca De Benítez'
$ws.Range("B75").Value = 'Cutzamala De Pinzón'
$ws.Range("B78").Value = 'Tixtla De Guerrero'
$ws.Range("B79").Value = 'Tlapa De Comonfort'
$ws.Range("B91").Value = 'Huejuquilla El Alto'
$ws.Range("B110").Value = 'Acatlán De Pérez Figueroa'
$ws.Range("B111").Value = 'Chalcatongo De Hidalgo'
$ws.Range("B113").Value = 'Coicoyán De Las Flores'
$ws.Range("B116").Value = 'Heroica Ciudad De Huajuapan De León'
$ws.Range("B117").Value = 'Ixtlán De Juárez'
$ws.Range("B118").Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Range("B122").Value = 'Miahuatlán De Porfirio Díaz'
$ws.Range("B123").Value = 'Oaxaca De Juárez'
$ws.Range("B124").Value = 'Putla Villa De Guerrero'
$ws.Range("B138").Value = 'Santo Domingo De Morelos'
$ws.Range("B142").Value = 'Tlacolula De Matamoros'
$ws.Range("B143").Value = 'Villa De Chilapa De Díaz'
$ws.Range("B144").Value = 'Zimatlán De Álvarez'
$ws.Range("B151").Value = 'Cuayuca De Andrade'
$ws.Range("B155").Value = 'Izúcar De Matamoros'
$ws.Range("B171").Value = 'Tepanco De López'
$ws.Range("B173").Value = 'Tlacotepec De Benito Juárez'
$ws.Range("B180").Value = 'Jalpan De Serra'
$ws.Range("B181").Value = 'Pinal De Amoles'
$ws.Range("B202").Value = 'Amatlán De Los Reyes'
$ws.Range("B206").Value = 'Cosamaloapan De Carpio'
$ws.Range("B209").Value = 'Martínez De La Torre'
$ws.Range("B212").Value = 'Paso Del Macho'
$ws.Range("B227").Value = 'Nochistlán De Mejía'

# Remove the footer/source notes rows (232-236)
$ws.Rows("232:236").Delete()
